$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - F1 label text stays "pvalues" (already existing shared string, reused)
$ws.Range("F1").Value = "pvalues"

# Column F now holds actual numeric p-values instead of the rounded text strings
$ws.Range("F2").Value = 0.20530027175263471
$ws.Range("F3").Value = 0.11858268365069501
$ws.Range("F4").Value = 0.2318148727428361
$ws.Range("F5").Value = 0.4795222310408439
$ws.Range("F6").Value = 0.97524781505563041

# Selection changed to whole column A
$ws.Range("A1:A1048576").Select()
